# Add logging (logstash) configuration rows to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A5").Value = "logstash_host"
$ws.Range("B5").Value = "10.2.100.56"
$ws.Range("A6").Value = "logstash_port"
$ws.Range("B6").Value = 5959

$ws.Range("B7").Select()
